$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.937.57'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.552.56'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').Value = "'206.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = "'21.56"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D10').Value = "'0.0583"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.773.20'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '1.551.16'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').Value = "'0.515"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '26.926.87'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = "'61.72"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').Value = "'213.97"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = "'7.25"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = "'1.00"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'4.03"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').Value = "'9.18"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('D25').Value = "'152.89"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  +2.29%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').Value = "'1.01"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('D30').Value = "'0.0460"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('D33').Value = '1.376.97'
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').Value = "'2.97"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('E35').Value = '  +3.30%  '
$ws.Range('D36').Value = "'0.970"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.82%  '
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').Value = "'0.808"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('D44').Value = "'2.25"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.36%  '
$ws.Range('D45').Value = "'63.72"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('D47').Value = '1.686.91'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').Value = "'86.23"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').Value = "'0.0509"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').Value = "'0.0953"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.92%  '
$ws.Range('E51').Value = '  +0.42%  '
